$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("INCIO SANCHEZ PAOLA KATHERINE", 25),
    @("GUEVARA IDROGO DENNIS PERCY", 25),
    @("TANTALEAN BUSTAMANTE ESTALIN YOEL", 24),
    @("HUAYHUA VALDIVIA LUZ EXMILDA", 23),
    @("CAMPOS PEREZ YOVERLY", 22),
    @("DELGADO VASQUEZ FLOR MAGALY", 20),
    @("MEDINA TAPIA ANA YULI", 20),
    @("LINARES PEREZ YANASELY", 20),
    @("PEREZ LINARES TATHIANA", 19),
    @("LOZADA ROJAS LUZ ELENA", 19),
    @("CHAVEZ VILLANUEVA SILVIA JANETH", 18),
    @("MONDRAGON HERNANDEZ WILMER JUNIOR", 18),
    @("SOTO LOZANO LUZDINA", 18),
    @("VASQUEZ SILVA ALOIS ADOLF", 17)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
